# 9.4.1.1 metadata sheet update
# - Organization / contact person / website (with hyperlink) refreshed
# - Active cell selection moved to H2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2. Информация об организации --------------------------------------------
$ws.Range("B6").Value = "Национальный статистический комитет КР (Управление цифрового развития и статистики устойчивого развития )"
$ws.Range("B7").Value = "Мамбеталиев Т.А."

# Сайт организации: new address + live hyperlink
$ws.Range("B10").Value = "www.stat.gov.kg "
$ws.Hyperlinks.Add($ws.Range("B10"), "http://www.stat.gov.kg/")

# Restore the active selection shown in the saved workbook (was A2 -> H2)
$ws.Range("H2").Select()

# Window / absolute-path metadata (best effort; harmless if unsupported by
# this headless host since these reflect desktop window chrome, not data)
try { $excel.ActiveWindow.Width = 28800 } catch {}
try {
    $wb.Application.ActiveWindow.Caption = "R:\Метаданные  ЦУР в Excel для Платформы\Метаданные на русском\Национальные\"
} catch {}
